$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the content of row 9 and row 10 for columns A, B, E, F, G, H, Q, R, S
# (these are the columns whose values differ between the two rows; the
# remaining columns already hold identical data on both rows).

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "S")

foreach ($col in $cols) {
    $addr9  = $col + "9"
    $addr10 = $col + "10"
    $v9  = $ws.Range($addr9).Value2
    $v10 = $ws.Range($addr10).Value2
    $ws.Range($addr9).Value2  = $v10
    $ws.Range($addr10).Value2 = $v9
}
